$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.185.65"
$ws.Range("E2").Value = "  +10.39%  "

# Row 3
$ws.Range("D3").Value = "1.680.00"
$ws.Range("E3").Value = "  +6.60%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.93"
$ws.Range("E5").Value = "  +7.34%  "

# Row 6
$ws.Range("E6").Value = "  +1.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3738"
$ws.Range("E7").Value = "  +0.96%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3444"
$ws.Range("E8").Value = "  +4.54%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.34"
$ws.Range("E9").Value = "  +12.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.192"
$ws.Range("E10").Value = "  +3.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07314"
$ws.Range("E11").Value = "  +3.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.43"
$ws.Range("E13").Value = "  +0.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.117"
$ws.Range("E14").Value = "  +4.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.779"
$ws.Range("E15").Value = "  +3.47%  "

# Row 16
$ws.Range("D16").Value = "1.675.76"
$ws.Range("E16").Value = "  +6.42%  "

# Row 17
$ws.Range("E17").Value = "  +2.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06733"
$ws.Range("E19").Value = "  +5.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.80"
$ws.Range("E20").Value = "  +7.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.123"
$ws.Range("E22").Value = "  +3.92%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.01"
$ws.Range("E23").Value = "  +2.70%  "

# Row 24
$ws.Range("D24").Value = "24.131.64"
$ws.Range("E24").Value = "  +10.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.418"
$ws.Range("E25").Value = "  +2.17%  "

# Row 26
$ws.Range("E26").Value = "  -9.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.673"
$ws.Range("E27").Value = "  +9.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.31"
$ws.Range("E28").Value = "  +2.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.62"
$ws.Range("E29").Value = "  +4.49%  "

# Row 30
$ws.Range("D30").Value = "1.862.57"
$ws.Range("E30").Value = "  +6.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.14"
$ws.Range("E31").Value = "  +5.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.411"
$ws.Range("E32").Value = "  +15.51%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.054"
$ws.Range("E33").Value = "  -2.86%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9801"
$ws.Range("E34").Value = "  +5.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.756"
$ws.Range("E35").Value = "  +8.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08460"
$ws.Range("E36").Value = "  +2.71%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.36"
$ws.Range("E37").Value = "  +3.87%  "

# Row 38
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.983"
$ws.Range("E38").Value = "  +3.58%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06442"
$ws.Range("E39").Value = "  +4.17%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.369"
$ws.Range("E40").Value = "  +2.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02350"
$ws.Range("E41").Value = "  +6.77%  "

# Row 42
$ws.Range("E42").Value = "  +2.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2126"
$ws.Range("E43").Value = "  +5.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6192"
$ws.Range("E44").Value = "  +5.90%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.807"
$ws.Range("E46").Value = "  +4.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.22"
$ws.Range("E47").Value = "  +1.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5983"
$ws.Range("E48").Value = "  +5.28%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043"
$ws.Range("E49").Value = "  +6.34%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.88"
$ws.Range("E50").Value = "  +0.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07162"
$ws.Range("E51").Value = "  +5.29%  "
